$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 136, shifting existing rows 136-224 down to 137-225.
$ws.Rows(136).Insert()

# Populate the newly inserted row 136 with the new record.
$ws.Range("A136").Value = 11
$ws.Range("B136").Value = "Vega Monumental Concepción"
$ws.Range("C136").Value = "Bíobío"
$ws.Range("D136").Value = 44813
$ws.Range("E136").Value = 8
$ws.Range("F136").Value = 100112040
$ws.Range("G136").Value = "Cilantro"
$ws.Range("H136").Value = "Sin especificar"
$ws.Range("I136").Value = "Primera"
$ws.Range("J136").Value = 200
$ws.Range("K136").Value = 5000
$ws.Range("L136").Value = 5500
$ws.Range("M136").Value = 5250
$ws.Range("N136").Value = "$/caja 36 atados"
$ws.Range("O136").Value = "Región Metropolitana"
$ws.Range("P136").Value = 146
$ws.Range("Q136").Value = 36
$ws.Range("R136").Value = "Hortaliza"
